$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new record to the cage database table (row 34)
$ws.Range("A34").Value = "453A"
$ws.Range("B34").Value = 3
$ws.Range("C34").Value = 34
$ws.Range("D34").Value = 22
$ws.Range("E34").Value = "Wood"
